# Add another author (He 2012, column E) to the author-modelling-comparison
# worksheet, and tweak two row labels in column A.
#
# Column E ("He 2012") already had a header cell (E1) but every data row in
# that column was blank - this fills those rows in. Two existing labels in
# column A are also reworded. A reviewer note is added at E14 using the
# built-in "Bad" (red) cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reword two existing row labels in column A -----------------------
$ws.Range("A5").Value() = "Processing of measured variable"
$ws.Range("A7").Value() = "Heath state- Health indicator mapping, fk"
$ws.Range("A8").Value() = "Health indicator - Measured variable mapping, hk"

# --- Fill in the new author's column (E = He 2012) ---------------------
$ws.Range("E2").Value()  = "Spiral bevel gear degradation (pitting)"
$ws.Range("E5").Value()  = "One dimensional transition function using whitening transform"
$ws.Range("E11").Value() = "Particle Filter with l-step ahead estimator"
$ws.Range("E10").Value() = "N/A This was buildt on data (ARIMA)"
$ws.Range("E7").Value()  = "Direct"
$ws.Range("E3").Value()  = "Oil debris"
$ws.Range("E4").Value()  = "Oil debris, Acceleration"
$ws.Range("E8").Value()  = " Data driven Double exponential smoothing model"

# --- Reviewer note, flagged with the built-in "Bad" style --------------
$ws.Range("E14").Value() = "Check if this is summarized in lit review"
$ws.Range("E14").Style = "Bad"

# --- Restore the selection left by the editor ---------------------------
[void]$ws.Range("F18").Select()
